# The deck ships two DrawingML theme parts: ppt/theme/theme1.xml ("Office
# Theme" colours) and ppt/theme/theme2.xml ("Integral" colours, currently
# wired to the one slide master / the whole design). The edit swaps the two
# themes' contents: the design that drives the slide master switches from
# the "Integral" palette to the standard "Office Theme" palette (font
# scheme and format scheme are identical between the two theme parts, so
# only the 12 colour-scheme entries actually change).
#
# PowerPoint's Theme object exposes exactly those 12 slots via
# ThemeColorScheme(1..12) in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB is read/written as a VBA colour long (0xBBGGRR), i.e. byte-reversed
# relative to the usual #RRGGBB hex notation.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0x000000   # dk1      -> 000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink -> 954F72
